$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.022.09"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "2.216.48"
$ws.Range("E3").Value = "  -1.43%  "
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("D5").Value = "'243.06"
$ws.Range("E5").Value = "  -1.72%  "
$ws.Range("D6").Value = "'0.626"
$ws.Range("E6").Value = "  +0.91%  "
$ws.Range("D7").Value = "'73.43"
$ws.Range("E7").Value = "  -1.73%  "
$ws.Range("D9").Value = "'0.611"
$ws.Range("E9").Value = "  -1.12%  "
$ws.Range("D10").Value = "'43.60"
$ws.Range("E10").Value = "  +5.34%  "
$ws.Range("D11").Value = "'0.0957"
$ws.Range("E11").Value = "  +1.54%  "
$ws.Range("D12").Value = "'7.10"
$ws.Range("E12").Value = "  +0.67%  "
$ws.Range("E13").Value = "  +1.27%  "
$ws.Range("D14").Value = "2.543.17"
$ws.Range("E14").Value = "  -1.50%  "
$ws.Range("D15").Value = "'14.23"
$ws.Range("E15").Value = "  -2.03%  "
$ws.Range("E16").Value = "  -1.17%  "
$ws.Range("D17").Value = "2.214.35"
$ws.Range("E17").Value = "  -1.72%  "
$ws.Range("D18").Value = "41.849.38"
$ws.Range("E18").Value = "  -0.52%  "
$ws.Range("E19").Value = "  +12.44%  "
$ws.Range("D20").Value = "'6.17"
$ws.Range("E20").Value = "  +0.38%  "
$ws.Range("D21").Value = "'72.31"
$ws.Range("E21").Value = "  +0.56%  "
$ws.Range("D22").Value = "'10.36"
$ws.Range("E22").Value = "  +34.46%  "
$ws.Range("D23").Value = "'229.20"
$ws.Range("E23").Value = "  -1.18%  "
$ws.Range("D24").Value = "'2.11"
$ws.Range("E24").Value = "  -8.39%  "
$ws.Range("D25").Value = "'11.55"
$ws.Range("E25").Value = "  +4.33%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("E27").Value = "  +1.20%  "
$ws.Range("D28").Value = "'2.27"
$ws.Range("E28").Value = "  -1.26%  "
$ws.Range("E29").Value = "  +5.88%  "
$ws.Range("D30").Value = "'166.66"
$ws.Range("E30").Value = "  -1.79%  "
$ws.Range("D31").Value = "'20.58"
$ws.Range("E31").Value = "  -0.43%  "
$ws.Range("D32").Value = "'5.58"
$ws.Range("E32").Value = "  +13.97%  "
$ws.Range("D33").Value = "'0.0795"
$ws.Range("E33").Value = "  -3.52%  "
$ws.Range("D34").Value = "'0.124"
$ws.Range("E34").Value = "  -0.27%  "
$ws.Range("D35").Value = "'29.34"
$ws.Range("E35").Value = "  -2.45%  "
$ws.Range("E36").Value = "  -4.22%  "
$ws.Range("E37").Value = "  -4.78%  "
$ws.Range("E38").Value = "  +0.08%  "
$ws.Range("D39").Value = "'12.94"
$ws.Range("E39").Value = "  -4.75%  "
$ws.Range("D40").Value = "'2.13"
$ws.Range("E40").Value = "  -2.61%  "
$ws.Range("D41").Value = "'64.46"
$ws.Range("E41").Value = "  +3.36%  "
$ws.Range("D42").Value = "'5.64"
$ws.Range("E42").Value = "  -2.69%  "
$ws.Range("E43").Value = "  -1.57%  "
$ws.Range("E44").Value = "  +0.28%  "
$ws.Range("D45").Value = "'103.96"
$ws.Range("E45").Value = "  -4.53%  "
$ws.Range("E46").Value = "  +0.69%  "
$ws.Range("E47").Value = "  +4.69%  "
$ws.Range("E48").Value = "  -0.68%  "
$ws.Range("E49").Value = "  +0.07%  "
$ws.Range("D50").Value = "'2.72"
$ws.Range("E50").Value = "  +0.71%  "
$ws.Range("D51").Value = "2.423.74"
$ws.Range("E51").Value = "  -1.31%  "
